# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled update).
# Numeric-looking text values are written with a leading apostrophe so Excel
# keeps them as literal text (preserving trailing zeros such as "1.00"/"308.48")
# instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '44.421.93'
$ws.Range('E2').Value = '  +1.15%  '

# Row 3
$ws.Range('D3').Value = '2.248.73'
$ws.Range('E3').Value = '  +0.86%  '

# Row 4
$ws.Range('E4').Value = '  +0.34%  '

# Row 5
$ws.Range('D5').Value = '''308.48'
$ws.Range('E5').Value = '  +1.69%  '

# Row 6
$ws.Range('D6').Value = '''94.93'
$ws.Range('E6').Value = '  +1.28%  '

# Row 7
$ws.Range('D7').Value = '''0.572'
$ws.Range('E7').Value = '  +1.12%  '

# Row 8
$ws.Range('D8').Value = '''1.00'
$ws.Range('E8').Value = '  +0.19%  '

# Row 9
$ws.Range('D9').Value = '''0.527'
$ws.Range('E9').Value = '  +2.12%  '

# Row 10
$ws.Range('D10').Value = '''35.50'
$ws.Range('E10').Value = '  +4.47%  '

# Row 11
$ws.Range('E11').Value = '  +2.28%  '

# Row 12
$ws.Range('E12').Value = '  +2.78%  '

# Row 13
$ws.Range('D13').Value = '''0.105'
$ws.Range('E13').Value = '  +2.01%  '

# Row 14
$ws.Range('D14').Value = '2.404.80'
$ws.Range('E14').Value = '  +5.87%  '

# Row 15
$ws.Range('D15').Value = '''0.841'
$ws.Range('E15').Value = '  +4.02%  '

# Row 16
$ws.Range('D16').Value = '''13.69'
$ws.Range('E16').Value = '  +2.05%  '

# Row 17
$ws.Range('D17').Value = '44.183.70'
$ws.Range('E17').Value = '  +1.18%  '

# Row 18
$ws.Range('D18').Value = '0.0₃0968'
$ws.Range('E18').Value = '  +2.04%  '

# Row 19
$ws.Range('D19').Value = '''12.31'
$ws.Range('E19').Value = '  +1.81%  '

# Row 20
$ws.Range('E20').Value = '  +5.24%  '

# Row 21
$ws.Range('E21').Value = '  +3.21%  '

# Row 22
$ws.Range('D22').Value = '''3.23'
$ws.Range('E22').Value = '  +11.82%  '

# Row 23
$ws.Range('D23').Value = '''237.83'
$ws.Range('E23').Value = '  +1.62%  '

# Row 24
$ws.Range('D24').Value = '''2.01'
$ws.Range('E24').Value = '  +5.77%  '

# Row 25
$ws.Range('E25').Value = '  -0.04%  '

# Row 27
$ws.Range('D27').Value = '''38.29'
$ws.Range('E27').Value = '  +6.77%  '

# Row 28
$ws.Range('D28').Value = '''9.87'
$ws.Range('E28').Value = '  +1.28%  '

# Row 29
$ws.Range('D29').Value = '''6.01'
$ws.Range('E29').Value = '  +1.98%  '

# Row 30
$ws.Range('D30').Value = '''20.14'
$ws.Range('E30').Value = '  +2.18%  '

# Row 31
$ws.Range('D31').Value = '''154.27'
$ws.Range('E31').Value = '  +2.47%  '

# Row 32
$ws.Range('D32').Value = '''0.0801'
$ws.Range('E32').Value = '  +0.73%  '

# Row 33
$ws.Range('E33').Value = '  +1.13%  '

# Row 34
$ws.Range('D34').Value = '''3.11'
$ws.Range('E34').Value = '  -2.26%  '

# Row 35
$ws.Range('E35').Value = '  +2.52%  '

# Row 36
$ws.Range('E36').Value = '  +3.74%  '

# Row 37
$ws.Range('E37').Value = '  +4.90%  '

# Row 38
$ws.Range('D38').Value = '''3.48'
$ws.Range('E38').Value = '  +7.43%  '

# Row 39
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Value = '''3.87'
$ws.Range('E39').Value = '  +2.93%  '

# Row 40
$ws.Range('B40').Value = 'Celestia'
$ws.Range('C40').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D40').Value = '''14.60'
$ws.Range('E40').Value = '  +1.41%  '

# Row 41
$ws.Range('E41').Value = '  +3.16%  '

# Row 42
$ws.Range('E42').Value = '  +0.32%  '

# Row 43
$ws.Range('D43').Value = '1.751.58'
$ws.Range('E43').Value = '  +1.31%  '

# Row 44
$ws.Range('E44').Value = '  +5.89%  '

# Row 45
$ws.Range('D45').Value = '''81.14'
$ws.Range('E45').Value = '  -2.54%  '

# Row 46
$ws.Range('D46').Value = '''71.24'
$ws.Range('E46').Value = '  +5.63%  '

# Row 47
$ws.Range('D47').Value = '''100.12'
$ws.Range('E47').Value = '  +1.08%  '

# Row 48
$ws.Range('E48').Value = '  +0.93%  '

# Row 49
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D49').Value = '''55.79'
$ws.Range('E49').Value = '  +4.41%  '

# Row 50
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').Value = '''1.60'
$ws.Range('E50').Value = '  +7.72%  '

# Row 51
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '''14.68'
$ws.Range('E51').Value = '  +0.93%  '
